$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1424070089188945
$ws.Range("D2").Value = 0.02022219814269377
$ws.Range("E2").Value = 0.1300736245576957
$ws.Range("F2").Value = 0.8825687668108486
$ws.Range("G2").Value = 0.7367339035159262
$ws.Range("H2").Value = 0.7828929359033054
$ws.Range("K2").Value = 0.6505724682276366
$ws.Range("M2").Value = 0.293988279558711
$ws.Range("N2").Value = 1.558108670163225

# Row 3
$ws.Range("B3").Value = 0.1329457677670831
$ws.Range("D3").Value = 0.01992685482388623
$ws.Range("E3").Value = 0.1227349444362389
$ws.Range("F3").Value = 0.8667623640130131
$ws.Range("G3").Value = 0.7210222625194405
$ws.Range("H3").Value = 0.7806078034431749
$ws.Range("K3").Value = 0.5680691593517508
$ws.Range("M3").Value = 0.2630141578407716
$ws.Range("N3").Value = 1.576700559340486

# Row 4
$ws.Range("B4").Value = 0.1272075476655914
$ws.Range("D4").Value = 0.01974403680389969
$ws.Range("E4").Value = 0.1183249178463512
$ws.Range("F4").Value = 0.857673231157662
$ws.Range("G4").Value = 0.7119386218113988
$ws.Range("H4").Value = 0.779648591287085
$ws.Range("K4").Value = 0.5173954926115414
$ws.Range("M4").Value = 0.2440992115783516
$ws.Range("N4").Value = 1.588707123821953

# Row 5
$ws.Range("B5").Value = 0.1248871700813083
$ws.Range("D5").Value = 0.019669173936812
$ws.Range("E5").Value = 0.1165516861483269
$ws.Range("F5").Value = 0.8541238063040595
$ws.Range("G5").Value = 0.7083780416049308
$ws.Range("H5").Value = 0.7793691639139411
$ws.Range("K5").Value = 0.4967414191117712
$ws.Range("M5").Value = 0.2364168823468233
$ws.Range("N5").Value = 1.593748372564473

# Row 6
$ws.Range("B6").Value = 0.1245029639070339
$ws.Range("D6").Value = 0.01965672132991259
$ws.Range("E6").Value = 0.1162586783321373
$ws.Range("F6").Value = 0.8535437437588484
$ws.Range("G6").Value = 0.7077953135808741
$ws.Range("H6").Value = 0.7793294925984782
$ws.Range("K6").Value = 0.4933115746945589
$ws.Range("M6").Value = 0.2351427786097346
$ws.Range("N6").Value = 1.594594431248787

# Row 7
$ws.Range("B7").Value = 0.1271761812061243
$ws.Range("D7").Value = 0.01974302863602873
$ws.Range("E7").Value = 0.1183009070073595
$ws.Range("F7").Value = 0.8576247375171988
$ws.Range("G7").Value = 0.7118900321448223
$ws.Range("H7").Value = 0.7796443717468691
$ws.Range("K7").Value = 0.5171169613772406
$ws.Range("M7").Value = 0.2439955014866655
$ws.Range("N7").Value = 1.588774510899283

# Row 8
$ws.Range("B8").Value = 0.1391301105050786
$ws.Range("D8").Value = 0.02012067484233881
$ws.Range("E8").Value = 0.1275231866572639
$ws.Range("F8").Value = 0.876990580579303
$ws.Range("G8").Value = 0.7311992414856121
$ws.Range("H8").Value = 0.7820128007097225
$ws.Range("K8").Value = 0.6221286577903982
$ws.Range("M8").Value = 0.283286769017316
$ws.Range("N8").Value = 1.564396386718673

# Row 9
$ws.Range("B9").Value = 0.163130623946941
$ws.Range("D9").Value = 0.02084921084797386
$ws.Range("E9").Value = 0.1463811732572324
$ws.Range("F9").Value = 0.9198784120681012
$ws.Range("G9").Value = 0.7735641526034556
$ws.Range("H9").Value = 0.7901875455598315
$ws.Range("K9").Value = 0.8279394793902668
$ws.Range("M9").Value = 0.361174317762277
$ws.Range("N9").Value = 1.521286940748958

# Row 10
$ws.Range("B10").Value = 0.1811002513312587
$ws.Range("D10").Value = 0.02137677097021751
$ws.Range("E10").Value = 0.1607252660449348
$ws.Range("F10").Value = 0.9544192688656921
$ws.Range("G10").Value = 0.8074779877264007
$ws.Range("H10").Value = 0.7983592453472568
$ws.Range("K10").Value = 0.9791101120300993
$ws.Range("M10").Value = 0.4189395078033158
$ws.Range("N10").Value = 1.492484555234899

# Row 11
$ws.Range("B11").Value = 0.1893473693060486
$ws.Range("D11").Value = 0.02161503062983883
$ws.Range("E11").Value = 0.1673609195727153
$ws.Range("F11").Value = 0.9707992878661003
$ws.Range("G11").Value = 0.8235217794801031
$ws.Range("H11").Value = 0.8025501154950518
$ws.Range("K11").Value = 1.047881180035461
$ws.Range("M11").Value = 0.4453428876159222
$ws.Range("N11").Value = 1.480006419416039

# Row 12
$ws.Range("B12").Value = 0.1924806740885145
$ws.Range("D12").Value = 0.02170499792265801
$ws.Range("E12").Value = 0.1698898417110044
$ws.Range("F12").Value = 0.9770984983490649
$ws.Range("G12").Value = 0.8296864853126635
$ws.Range("H12").Value = 0.8042053993576701
$ws.Range("K12").Value = 1.073923711404632
$ws.Range("M12").Value = 0.4553596770339112
$ws.Range("N12").Value = 1.475371212329016

# Row 13
$ws.Range("B13").Value = 0.1918054050627518
$ws.Range("D13").Value = 0.02168563338321405
$ws.Range("E13").Value = 0.1693444708884471
$ws.Range("F13").Value = 0.9757375527632064
$ws.Range("G13").Value = 0.8283548247732426
$ws.Range("H13").Value = 0.8038458638261261
$ws.Range("K13").Value = 1.068314967237541
$ws.Range("M13").Value = 0.4532015585219256
$ws.Range("N13").Value = 1.476365480963242

# Row 14
$ws.Range("B14").Value = 0.1896049425813544
$ws.Range("D14").Value = 0.02162243747639891
$ws.Range("E14").Value = 0.1675686504641547
$ws.Range("F14").Value = 0.9713155917206677
$ws.Range("G14").Value = 0.8240271606484555
$ws.Range("H14").Value = 0.8026849268598824
$ws.Range("K14").Value = 1.050023706204456
$ws.Range("M14").Value = 0.4461666042150654
$ws.Range("N14").Value = 1.479623273226302

# Row 15
$ws.Range("B15").Value = 0.1882584328288033
$ws.Range("D15").Value = 0.02158369453468367
$ws.Range("E15").Value = 0.1664830194037208
$ws.Range("F15").Value = 0.9686195890547538
$ws.Range("G15").Value = 0.8213879866214313
$ws.Range("H15").Value = 0.8019827191043873
$ws.Range("K15").Value = 1.038819836111145
$ws.Range("M15").Value = 0.4418598978634947
$ws.Range("N15").Value = 1.481630491608133

# Row 16
$ws.Range("B16").Value = 0.1805627320904222
$ws.Range("D16").Value = 0.02136116466510529
$ws.Range("E16").Value = 0.160293855366767
$ws.Range("F16").Value = 0.9533622654927285
$ws.Range("G16").Value = 0.8064419443149404
$ws.Range("H16").Value = 0.7980949076818433
$ws.Range("K16").Value = 0.9746158186366358
$ws.Range("M16").Value = 0.4172165421602898
$ws.Range("N16").Value = 1.493312611184965

# Row 17
$ws.Range("B17").Value = 0.175860171079492
$ws.Range("D17").Value = 0.02122420114663015
$ws.Range("E17").Value = 0.1565254787180734
$ws.Range("F17").Value = 0.9441736411756665
$ws.Range("G17").Value = 0.7974313272925144
$ws.Range("H17").Value = 0.795831283613353
$ws.Range("K17").Value = 0.9352294596108948
$ws.Range("M17").Value = 0.402131068822186
$ws.Range("N17").Value = 1.50063925236423

# Row 18
$ws.Range("B18").Value = 0.1731622287285859
$ws.Range("D18").Value = 0.0211452609543592
$ws.Range("E18").Value = 0.1543683888606537
$ws.Range("F18").Value = 0.9389513578339859
$ws.Range("G18").Value = 0.792306620467599
$ws.Range("H18").Value = 0.7945738649945326
$ws.Range("K18").Value = 0.9125757286068108
$ws.Range("M18").Value = 0.3934661158188106
$ws.Range("N18").Value = 1.504912082460589

# Row 19
$ws.Range("B19").Value = 0.1722499321206641
$ws.Range("D19").Value = 0.0211185055073102
$ws.Range("E19").Value = 0.1536398102718053
$ws.Range("F19").Value = 0.9371939478461258
$ws.Range("G19").Value = 0.7905814177695447
$ws.Range("H19").Value = 0.7941557714927399
$ws.Range("K19").Value = 0.9049056088641407
$ws.Range("M19").Value = 0.3905343304937219
$ws.Range("N19").Value = 1.506368870981415

# Row 20
$ws.Range("B20").Value = 0.1763600592118877
$ws.Range("D20").Value = 0.02123879801147055
$ws.Range("E20").Value = 0.1569255525965829
$ws.Range("F20").Value = 0.9451452857470031
$ws.Range("G20").Value = 0.7983845192687795
$ws.Range("H20").Value = 0.7960676372177033
$ws.Range("K20").Value = 0.9394221776422
$ws.Range("M20").Value = 0.4037357183381545
$ws.Range("N20").Value = 1.499853237878536

# Row 21
$ws.Range("B21").Value = 0.1902509933722314
$ws.Range("D21").Value = 0.02164100666709956
$ws.Range("E21").Value = 0.1680898117319671
$ws.Range("F21").Value = 0.9726118070968255
$ws.Range("G21").Value = 0.825295872474868
$ws.Range("H21").Value = 0.8030240671468221
$ws.Range("K21").Value = 1.05539628135989
$ws.Range("M21").Value = 0.4482324387691392
$ws.Range("N21").Value = 1.478663936446967

# Row 22
$ws.Range("B22").Value = 0.1993894914255208
$ws.Range("D22").Value = 0.02190237496097325
$ws.Range("E22").Value = 0.1754805234478454
$ws.Range("F22").Value = 0.9911252056476485
$ws.Range("G22").Value = 0.8434046400481918
$ws.Range("H22").Value = 0.8079686506862345
$ws.Range("K22").Value = 1.131194526586512
$ws.Range("M22").Value = 0.4774210355971888
$ws.Range("N22").Value = 1.465340156105107

# Row 23
$ws.Range("B23").Value = 0.1945066640685127
$ws.Range("D23").Value = 0.02176301736358255
$ws.Range("E23").Value = 0.1715272533660013
$ws.Range("F23").Value = 0.9811926242562095
$ws.Range("G23").Value = 0.8336917926269223
$ws.Range("H23").Value = 0.8052931354293378
$ws.Range("K23").Value = 1.090739343038365
$ws.Range("M23").Value = 0.4618325947167961
$ws.Range("N23").Value = 1.472403225305079

# Row 24
$ws.Range("B24").Value = 0.1761340422629161
$ws.Range("D24").Value = 0.02123219938373921
$ws.Range("E24").Value = 0.1567446499389504
$ws.Range("F24").Value = 0.9447058171916609
$ws.Range("G24").Value = 0.7979534079791506
$ws.Range("H24").Value = 0.7959606447979581
$ws.Range("K24").Value = 0.9375266809463199
$ws.Range("M24").Value = 0.4030102317586426
$ws.Range("N24").Value = 1.500208406121086

# Row 25
$ws.Range("B25").Value = 0.1565784565892727
$ws.Range("D25").Value = 0.02065344848825745
$ws.Range("E25").Value = 0.1411949849309266
$ws.Range("F25").Value = 0.9077464973700415
$ws.Range("G25").Value = 0.7616169422042418
$ws.Range("H25").Value = 0.7875967277125397
$ws.Range("K25").Value = 0.7722713536255128
$ws.Range("M25").Value = 0.293988279558711
$ws.Range("N25").Value = 1.532445708852578
